# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Ají" (Hortaliza) at the top of the
# data block (rows 57-58), pushing the existing rows 57-84 down to 59-86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 57:58 - everything currently at/after row 57
# (through row 84) shifts down to rows 59-86, carrying its formatting
# (including column D's date style) along with it.
$ws.Rows("57:58").Insert()

# New row 57: Ají "Inferno", calidad "Primera"
$ws.Range("A57").Value2 = 8
$ws.Range("B57").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C57").Value2 = "Coquimbo"
$ws.Range("D57").Value2 = 44438
$ws.Range("E57").Value2 = 4
$ws.Range("F57").Value2 = 100112021
$ws.Range("G57").Value2 = "Ají"
$ws.Range("H57").Value2 = "Inferno"
$ws.Range("I57").Value2 = "Primera"
$ws.Range("J57").Value2 = 600
$ws.Range("K57").Value2 = 36000
$ws.Range("L57").Value2 = 37000
$ws.Range("M57").Value2 = 36500
$ws.Range("N57").Value2 = "$/caja 12 kilos"
$ws.Range("O57").Value2 = "Región de Arica y Parinacota"
$ws.Range("P57").Value2 = 3042
$ws.Range("Q57").Value2 = 12
$ws.Range("R57").Value2 = "Hortaliza"

# New row 58: Ají "Inferno", calidad "Segunda"
$ws.Range("A58").Value2 = 8
$ws.Range("B58").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C58").Value2 = "Coquimbo"
$ws.Range("D58").Value2 = 44438
$ws.Range("E58").Value2 = 4
$ws.Range("F58").Value2 = 100112021
$ws.Range("G58").Value2 = "Ají"
$ws.Range("H58").Value2 = "Inferno"
$ws.Range("I58").Value2 = "Segunda"
$ws.Range("J58").Value2 = 400
$ws.Range("K58").Value2 = 30000
$ws.Range("L58").Value2 = 31000
$ws.Range("M58").Value2 = 30500
$ws.Range("N58").Value2 = "$/caja 12 kilos"
$ws.Range("O58").Value2 = "Región de Arica y Parinacota"
$ws.Range("P58").Value2 = 2542
$ws.Range("Q58").Value2 = 12
$ws.Range("R58").Value2 = "Hortaliza"
